$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking price strings
# are preserved as text (matching the original inlineStr cell type),
# then revert the style so no stray cell formatting is introduced.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.449.32'
$ws.Range("E2").Value = '  -3.68%  '
$ws.Range("D3").Value = '1.851.82'
$ws.Range("E3").Value = '  -5.41%  '
$ws.Range("E4").Value = '  -0.78%  '
$ws.Range("D5").Value = '320.81'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("D7").Value = '0.4467'
$ws.Range("E7").Value = '  -6.26%  '
$ws.Range("D8").Value = '0.3836'
$ws.Range("E8").Value = '  -5.71%  '
$ws.Range("D9").Value = '47.95'
$ws.Range("E9").Value = '  -10.29%  '
$ws.Range("D10").Value = '0.07820'
$ws.Range("E10").Value = '  -7.43%  '
$ws.Range("D11").Value = '1.015'
$ws.Range("E11").Value = '  -4.29%  '
$ws.Range("D12").Value = '21.37'
$ws.Range("E12").Value = '  -3.31%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.842.03'
$ws.Range("E13").Value = '  -6.33%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.843'
$ws.Range("E14").Value = '  -5.59%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.093'
$ws.Range("E15").Value = '  -7.01%  '
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '85.61'
$ws.Range("E17").Value = '  -4.11%  '
$ws.Range("D18").Value = '0.00001022'
$ws.Range("E18").Value = '  -4.77%  '
$ws.Range("D19").Value = '0.06493'
$ws.Range("E19").Value = '  -1.98%  '
$ws.Range("D20").Value = '17.02'
$ws.Range("E20").Value = '  -9.22%  '
$ws.Range("E21").Value = '  -0.82%  '
$ws.Range("D22").Value = '5.470'
$ws.Range("E22").Value = '  -6.12%  '
$ws.Range("D23").Value = '27.386.87'
$ws.Range("E23").Value = '  -3.94%  '
$ws.Range("D24").Value = '10.75'
$ws.Range("E24").Value = '  -6.99%  '
$ws.Range("D25").Value = '2.287'
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").Value = '2.008.62'
$ws.Range("E26").Value = '  -8.87%  '
$ws.Range("D27").Value = '150.87'
$ws.Range("E27").Value = '  -2.21%  '
$ws.Range("E28").Value = '  -4.43%  '
$ws.Range("D29").Value = '5.459'
$ws.Range("E29").Value = '  -8.91%  '
$ws.Range("D30").Value = '2.024'
$ws.Range("E30").Value = '  -6.50%  '
$ws.Range("D31").Value = '119.60'
$ws.Range("E31").Value = '  -3.41%  '
$ws.Range("D32").Value = '1.496'
$ws.Range("E32").Value = '  +2.83%  '
$ws.Range("D33").Value = '0.09336'
$ws.Range("E33").Value = '  -2.72%  '
$ws.Range("D34").Value = '0.9253'
$ws.Range("E34").Value = '  -6.63%  '
$ws.Range("D35").Value = '3.617'
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").Value = '5.209'
$ws.Range("E36").Value = '  -6.91%  '
$ws.Range("E37").Value = '  -2.97%  '
$ws.Range("D38").Value = '0.02215'
$ws.Range("E38").Value = '  -5.06%  '
$ws.Range("D39").Value = '0.05931'
$ws.Range("E39").Value = '  -4.75%  '
$ws.Range("D40").Value = '8.269'
$ws.Range("E40").Value = '  -6.18%  '
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("D42").Value = '0.5879'
$ws.Range("E42").Value = '  -5.71%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '0.1847'
$ws.Range("E43").Value = '  -4.08%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '10.23'
$ws.Range("E44").Value = '  -8.32%  '
$ws.Range("D45").Value = '1.278'
$ws.Range("E45").Value = '  -4.22%  '
$ws.Range("D46").Value = '0.5597'
$ws.Range("E46").Value = '  -6.30%  '
$ws.Range("D47").Value = '12.14'
$ws.Range("E47").Value = '  -7.06%  '
$ws.Range("D48").Value = '3.352'
$ws.Range("E48").Value = '  -1.44%  '
$ws.Range("D49").Value = '1.903'
$ws.Range("E49").Value = '  -7.57%  '
$ws.Range("D50").Value = '0.06841'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("E51").Value = '  -10.99%  '

$dataRange.Style = "Normal"

